# Update countries & provincias Spain
# Refreshes the COVID-19 country statistics (sheet "Pais") to the
# 20:22 snapshot, including the ranking re-sort (Barein overtakes
# Azerbaiyan; Ruanda overtakes Trinidad yTobago) and the timestamp
# footer text.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Pais")

# --- Timestamp footer (row 1) ---
$ws.Range("A1").Value = "Datos actualizados a 8 de Abril de 2020 a las 20:22"

# --- Estados Unidos (row 4) ---
$ws.Range("B4").Value = 418410
$ws.Range("C4").Value = 18075
$ws.Range("E4").Value = 381986
$ws.Range("G4").Value = 1399
$ws.Range("H4").Value = 14240

# --- Alemania (row 8) ---
$ws.Range("B8").Value = 110483
$ws.Range("C8").Value = 2820
$ws.Range("E8").Value = 72219
$ws.Range("G8").Value = 167
$ws.Range("H8").Value = 2183

# --- Canada (row 16) ---
$ws.Range("B16").Value = 19179
$ws.Range("C16").Value = 1282
$ws.Range("D16").Value = 4474
$ws.Range("E16").Value = 14278
$ws.Range("G16").Value = 46
$ws.Range("H16").Value = 427

# --- Brasil (row 17) ---
$ws.Range("B17").Value = 14347
$ws.Range("C17").Value = 313
$ws.Range("E17").Value = 13501
$ws.Range("G17").Value = 33
$ws.Range("H17").Value = 719

# --- Austria (row 19) ---
$ws.Range("B19").Value = 12930
$ws.Range("C19").Value = 291
$ws.Range("E19").Value = 8145

# --- Egipto (row 58) ---
$ws.Range("B58").Value = 1560
$ws.Range("C58").Value = 110
$ws.Range("D58").Value = 305
$ws.Range("E58").Value = 1152
$ws.Range("G58").Value = 9
$ws.Range("H58").Value = 103

# --- Barein overtakes Azerbaiyan: rows 72/73 swap ---
$ws.Range("A72").Value = "Barein"
$ws.Range("B72").Value = 823
$ws.Range("C72").Value = 12
$ws.Range("D72").Value = 477
$ws.Range("E72").Value = 341
$ws.Range("F72").Value = 3
$ws.Range("H72").Value = 5

$ws.Range("A73").Value = "Azerbaiyan"
$ws.Range("B73").Value = 822
$ws.Range("C73").Value = 105
$ws.Range("D73").Value = 63
$ws.Range("E73").Value = 751
$ws.Range("F73").Value = 23
$ws.Range("H73").Value = 8

# --- Ruanda overtakes Trinidad yTobago: rows 128/129 swap ---
$ws.Range("A128").Value = "Ruanda"
$ws.Range("B128").Value = 110
$ws.Range("C128").Value = 5
$ws.Range("D128").Value = 7
$ws.Range("E128").Value = 103
$ws.Range("H128").Value = 0

$ws.Range("A129").Value = "Trinidad yTobago"
$ws.Range("B129").Value = 107
$ws.Range("D129").Value = 1
$ws.Range("H129").Value = 8

# --- Zambia (row 149) ---
$ws.Range("F149").Value = 1
